$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Paragraph 1: "n=11, m=19*2=38" run gets marked dirty ---
$para1 = $tr.Paragraphs(1, 1)
$run1 = $para1.Runs(1, 1)
$run1.Text = "n=11, m=19*2=38"

# --- Paragraph 3: split the "Inactive links..." run in two, updating the edge list ---
$para3 = $tr.Paragraphs(3, 1)
$run3 = $para3.Runs(1, 1)
$run3.Text = "Inactive links after first TC iteration (k=2.0"
$run3.InsertAfter("): e1-3,e2-3,e2-4,e2-5,e3-9,e3-11,e9-11")

# --- Mark edge 3-9 (node "9" <-> node "3") as inactive by giving it a dotted line ---
$grp = $s.Shapes.Item(3)
$edge39 = $grp.GroupItems.Item(14)
$edge39.Line.DashStyle = 2
